$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "have"
$ws.Range("B6").Value = "có"

$ws.Range("A7").Value = "havaaa"
$ws.Range("B7").Value = "havaaa"

$ws.Range("A8").Value = "haw"
$ws.Range("B8").Value = "ồ"
